$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Dashboard"
#   - H3:H23 (PERIOD TO EXPIRE) drops by 8 days for every row.
#   - I3:I23 (LAST UPDATE) moves 8 days later (08-Sep-2025 -> 16-Sep-2025),
#     staying a literal "dd-MMM-yyyy" text string (not a real date), and
#     keeping each row's original cell style.
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item(1)
$months = @("Jan","Feb","Mar","Apr","May","Jun","Jul","Aug","Sep","Oct","Nov","Dec")

for ($r = 3; $r -le 23; $r++) {
    # PERIOD TO EXPIRE shrinks by 8.
    $periodCell = $wsTraining.Cells.Item($r, 8)
    $periodCell.Value = $periodCell.Value2 - 8

    # LAST UPDATE advances by 8 days, formatted the same way as the source data.
    $updateCell = $wsTraining.Cells.Item($r, 9)
    $parts = $updateCell.Value2.Split("-")
    $day = [int]$parts[0]
    $monthIndex = [array]::IndexOf($months, $parts[1]) + 1
    $year = [int]$parts[2]
    $newDate = (Get-Date -Year $year -Month $monthIndex -Day $day).AddDays(8)
    $newDateText = "{0:D2}-{1}-{2}" -f $newDate.Day, $months[$newDate.Month - 1], $newDate.Year

    # Column J on the same row never changes, so borrow its pristine style as
    # a format stamp -- this keeps the rewritten cell on its original style
    # (s="3" / s="4") instead of drifting to a new one just because we had to
    # force a text number format to stop "16-Sep-2025" being read as a date.
    $styleDonor = $wsTraining.Cells.Item($r, 10)
    $updateCell.NumberFormat = "@"
    $updateCell.Value = $newDateText
    $styleDonor.Copy()
    $updateCell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Exam Dashboard"
#   - E3:E8 (COMMENTS) text changes from "OK" to "date is valid".
#   - Column E gets wider (10 -> 15) to fit the longer text.
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item(2)

for ($r = 3; $r -le 8; $r++) {
    $wsExam.Cells.Item($r, 5).Value = "date is valid"
}

$wsExam.Columns.Item(5).ColumnWidth = 14.17

# ---------------------------------------------------------------------------
# Header styling on both sheets:
#   - The bold title cell (A1) loses its 14pt size and turns white.
#   - The bold header row (row 2, dark-blue fill) also turns white.
#   (Both end up sharing the same bold/white font.)
# ---------------------------------------------------------------------------
$sheetInfo = @{1 = 11; 2 = 7}
foreach ($idx in $sheetInfo.Keys) {
    $ws = $wb.Worksheets.Item($idx)
    $lastCol = $sheetInfo[$idx]

    $titleCell = $ws.Range("A1")
    $titleCell.Font.Size = 11
    $titleCell.Font.Color = 16777215

    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Color = 16777215
}
